$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.131.90"
$ws.Range("E2").Value = "  +0.46%  "

$ws.Range("D3").Value = "2.903.71"
$ws.Range("E3").Value = "  +3.40%  "

$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "352.00"
$ws.Range("E5").Value = "  -0.59%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "113.14"
$ws.Range("E6").Value = "  +0.89%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.556"
$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.620"
$ws.Range("E9").Value = "  -0.76%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.64"
$ws.Range("E10").Value = "  -1.84%  "

$ws.Range("E11").Value = "  +0.69%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0864"
$ws.Range("E12").Value = "  +3.03%  "

$ws.Range("E13").Value = "  -1.11%  "

$ws.Range("E14").Value = "  -1.17%  "

$ws.Range("D15").Value = "3.359.20"
$ws.Range("E15").Value = "  +3.46%  "

$ws.Range("D16").Value = "2.904.48"
$ws.Range("E16").Value = "  +4.05%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.984"
$ws.Range("E17").Value = "  +4.13%  "

$ws.Range("D18").Value = "52.231.92"
$ws.Range("E18").Value = "  +0.77%  "

$ws.Range("B19").Value = "ImmutableX"
$ws.Range("C19").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.32"
$ws.Range("E19").Value = "  +2.91%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.59"
$ws.Range("E20").Value = "  -0.56%  "

$ws.Range("E21").Value = "  +2.29%  "

$ws.Range("D22").Value = "0.0₃0975"
$ws.Range("E22").Value = "  +0.16%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.01"
$ws.Range("E23").Value = "  +0.91%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.91"
$ws.Range("E24").Value = "  +0.33%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.80"
$ws.Range("E25").Value = "  +1.21%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.183"
$ws.Range("E26").Value = "  +14.59%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.75"
$ws.Range("E27").Value = "  +2.18%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.10%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.63"
$ws.Range("E29").Value = "  +2.36%  "

$ws.Range("E30").Value = "  +16.07%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.68"
$ws.Range("E31").Value = "  +8.66%  "

$ws.Range("E32").Value = "  -4.61%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.27"
$ws.Range("E33").Value = "  -0.13%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.17"
$ws.Range("E34").Value = "  +11.28%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "53.33"
$ws.Range("E35").Value = "  +1.53%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0449"
$ws.Range("E36").Value = "  -0.86%  "

$ws.Range("E37").Value = "  -0.11%  "

$ws.Range("E38").Value = "  +4.27%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.82"
$ws.Range("E39").Value = "  -0.59%  "

$ws.Range("E40").Value = "  +1.25%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.73"
$ws.Range("E41").Value = "  +7.61%  "

$ws.Range("E42").Value = "  +1.32%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.92"
$ws.Range("E43").Value = "  +4.77%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "119.51"
$ws.Range("E44").Value = "  -0.65%  "

$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.19"
$ws.Range("E45").Value = "  -1.91%  "

$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.57"
$ws.Range("E46").Value = "  +4.51%  "

$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "2.167.89"
$ws.Range("E47").Value = "  +2.89%  "

$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.48"
$ws.Range("E48").Value = "  -0.05%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.263"
$ws.Range("E49").Value = "  +20.25%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0344"
$ws.Range("E50").Value = "  +10.51%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.954"
$ws.Range("E51").Value = "  -1.25%  "

